$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.220.49'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '1.860.70'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7129'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '237.91'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.08199'
$ws.Range('E8').Value = '  +10.79%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3045'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.21'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '1.848.11'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.176'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7089'
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').Value = '29.242.55'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000007942'
$ws.Range('E17').Value = '  +3.92%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.791'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('E19').Value = '  +2.30%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '237.34'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '2.109.53'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.430'
$ws.Range('E24').Value = '  -2.20%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '162.66'
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.959'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1459'
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.09'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.962'
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.486'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.404'
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.026'
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05224'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7085'
$ws.Range('E36').Value = '  +0.83%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.003'
$ws.Range('E37').Value = '  -2.80%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.677'
$ws.Range('E38').Value = '  +0.38%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01859'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.729'
$ws.Range('E40').Value = '  +2.10%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9240'
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('D42').Value = '1.140.23'
$ws.Range('E42').Value = '  +6.36%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.4285'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.901'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '70.22'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.9995'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '102.83'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.777'
$ws.Range('E48').Value = '  +2.06%  '
$ws.Range('D49').Value = '2.009.20'
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.219'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.954'
$ws.Range('E51').Value = '  -0.95%  '
